$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'22.158.04"
$ws.Range("D2").Style = 'Normal'
$ws.Range("E2").Value = '  -1.28%  '

$ws.Range("D3").Value = "'1.560.03"
$ws.Range("D3").Style = 'Normal'
$ws.Range("E3").Value = '  -0.62%  '

$ws.Range("D4").Value = "'0.9987"
$ws.Range("D4").Style = 'Normal'
$ws.Range("E4").Value = '  -0.30%  '

$ws.Range("D5").Value = "'0.9984"
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  -0.33%  '

$ws.Range("D6").Value = "'288.54"
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = '  -0.01%  '

$ws.Range("D7").Value = "'0.3801"
$ws.Range("D7").Style = 'Normal'
$ws.Range("E7").Value = '  +2.61%  '

$ws.Range("D8").Value = "'0.3299"
$ws.Range("D8").Style = 'Normal'
$ws.Range("E8").Value = '  -0.24%  '

$ws.Range("D9").Value = "'43.45"
$ws.Range("D9").Style = 'Normal'
$ws.Range("E9").Value = '  -10.01%  '

$ws.Range("D10").Value = "'1.147"
$ws.Range("D10").Style = 'Normal'
$ws.Range("E10").Value = '  +1.04%  '

$ws.Range("D11").Value = "'0.07387"
$ws.Range("D11").Style = 'Normal'
$ws.Range("E11").Value = '  -1.47%  '

$ws.Range("D12").Value = "'0.9989"
$ws.Range("D12").Style = 'Normal'
$ws.Range("E12").Value = '  -0.29%  '

$ws.Range("D13").Value = "'20.22"
$ws.Range("D13").Style = 'Normal'
$ws.Range("E13").Value = '  -2.38%  '

$ws.Range("D14").Value = "'5.842"
$ws.Range("D14").Style = 'Normal'
$ws.Range("E14").Value = '  -1.44%  '

$ws.Range("D15").Value = "'6.871"
$ws.Range("D15").Style = 'Normal'
$ws.Range("E15").Value = '  +0.02%  '

$ws.Range("D16").Value = "'1.567.28"
$ws.Range("D16").Style = 'Normal'
$ws.Range("E16").Value = '  -0.14%  '

$ws.Range("D17").Value = "'0.00001107"
$ws.Range("D17").Style = 'Normal'
$ws.Range("E17").Value = '  -0.99%  '

$ws.Range("B18").Value = 'Litecoin'
$ws.Range("C18").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D18").Value = "'86.18"
$ws.Range("D18").Style = 'Normal'
$ws.Range("E18").Value = '  -1.47%  '

$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").Value = "'0.06625"
$ws.Range("D19").Style = 'Normal'
$ws.Range("E19").Value = '  -1.87%  '

$ws.Range("D20").Value = "'6.421"
$ws.Range("D20").Style = 'Normal'
$ws.Range("E20").Value = '  +1.25%  '

$ws.Range("D21").Value = "'0.9980"
$ws.Range("D21").Style = 'Normal'
$ws.Range("E21").Value = '  -0.31%  '

$ws.Range("E22").Value = '  -2.71%  '

$ws.Range("D23").Value = "'11.72"
$ws.Range("D23").Style = 'Normal'
$ws.Range("E23").Value = '  -2.51%  '

$ws.Range("D24").Value = "'22.150.37"
$ws.Range("D24").Style = 'Normal'
$ws.Range("E24").Value = '  -1.31%  '

$ws.Range("D25").Value = "'2.311"
$ws.Range("D25").Style = 'Normal'

$ws.Range("D26").Value = "'2.537"
$ws.Range("D26").Style = 'Normal'
$ws.Range("E26").Value = '  -1.94%  '

$ws.Range("D27").Value = "'150.20"
$ws.Range("D27").Style = 'Normal'
$ws.Range("E27").Value = '  -2.24%  '

$ws.Range("D28").Value = "'19.18"
$ws.Range("D28").Style = 'Normal'
$ws.Range("E28").Value = '  -2.59%  '

$ws.Range("D29").Value = "'4.936"
$ws.Range("D29").Style = 'Normal'
$ws.Range("E29").Value = '  -1.76%  '

$ws.Range("D30").Value = "'121.96"
$ws.Range("D30").Style = 'Normal'
$ws.Range("E30").Value = '  -1.93%  '

$ws.Range("D31").Value = "'1.739.91"
$ws.Range("D31").Style = 'Normal'
$ws.Range("E31").Value = '  -0.29%  '

$ws.Range("D32").Value = "'1.080"
$ws.Range("D32").Style = 'Normal'
$ws.Range("E32").Value = '  +2.38%  '

$ws.Range("D33").Value = "'5.957"
$ws.Range("D33").Style = 'Normal'
$ws.Range("E33").Value = '  -2.50%  '

$ws.Range("D34").Value = "'1.849"
$ws.Range("D34").Style = 'Normal'
$ws.Range("E34").Value = '  -8.15%  '

$ws.Range("D35").Value = "'0.08264"
$ws.Range("D35").Style = 'Normal'
$ws.Range("E35").Value = '  -1.24%  '

$ws.Range("D36").Value = "'9.357"
$ws.Range("D36").Style = 'Normal'
$ws.Range("E36").Value = '  -4.39%  '

$ws.Range("D37").Value = "'0.02351"
$ws.Range("D37").Style = 'Normal'
$ws.Range("E37").Value = '  -4.71%  '

$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D38").Value = "'5.318"
$ws.Range("D38").Style = 'Normal'
$ws.Range("E38").Value = '  -0.24%  '

$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").Value = "'0.06269"
$ws.Range("D39").Style = 'Normal'
$ws.Range("E39").Value = '  -2.07%  '

$ws.Range("D40").Value = "'0.2172"
$ws.Range("D40").Style = 'Normal'
$ws.Range("E40").Value = '  -3.85%  '

$ws.Range("E41").Value = '  -2.53%  '

$ws.Range("E42").Value = '  -2.05%  '

$ws.Range("D43").Value = "'0.6089"
$ws.Range("D43").Style = 'Normal'
$ws.Range("E43").Value = '  -3.58%  '

$ws.Range("D44").Value = "'0.9982"
$ws.Range("D44").Style = 'Normal'
$ws.Range("E44").Value = '  -0.29%  '

$ws.Range("D45").Value = "'13.81"
$ws.Range("D45").Style = 'Normal'
$ws.Range("E45").Value = '  +0.16%  '

$ws.Range("D46").Value = "'3.743"
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = '  -0.73%  '

$ws.Range("D47").Value = "'0.5897"
$ws.Range("D47").Style = 'Normal'
$ws.Range("E47").Value = '  -4.10%  '

$ws.Range("D48").Value = "'1.997"
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = '  -3.05%  '

$ws.Range("D49").Value = "'122.33"
$ws.Range("D49").Style = 'Normal'
$ws.Range("E49").Value = '  -3.07%  '

$ws.Range("E50").Value = '  -2.81%  '

$ws.Range("D51").Value = "'0.07030"
$ws.Range("D51").Style = 'Normal'
$ws.Range("E51").Value = '  -2.59%  '
